$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.718.64'
$ws.Range('E2').Value = '  +3.05%  '
$ws.Range('D3').Value = '1.864.64'
$ws.Range('E3').Value = '  +2.91%  '
$ws.Range('D4').Value = '''1.034'
$ws.Range('E4').Value = '  +2.60%  '
$ws.Range('D5').Value = '''325.00'
$ws.Range('E5').Value = '  +3.45%  '
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('D7').Value = '''0.4411'
$ws.Range('E7').Value = '  +2.79%  '
$ws.Range('D8').Value = '''0.3801'
$ws.Range('E8').Value = '  +2.95%  '
$ws.Range('D9').Value = '''0.07462'
$ws.Range('E9').Value = '  +3.06%  '
$ws.Range('D10').Value = '''0.8846'
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('D11').Value = '''21.79'
$ws.Range('E11').Value = '  +2.41%  '
$ws.Range('D12').Value = '1.881.75'
$ws.Range('E12').Value = '  -8.03%  '
$ws.Range('D13').Value = '''5.554'
$ws.Range('E13').Value = '  +2.77%  '
$ws.Range('D14').Value = '''6.753'
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('D15').Value = '''0.07211'
$ws.Range('E15').Value = '  +3.78%  '
$ws.Range('D16').Value = '''83.83'
$ws.Range('E16').Value = '  +3.73%  '
$ws.Range('E17').Value = '  +3.02%  '
$ws.Range('D18').Value = '''0.000009109'
$ws.Range('E18').Value = '  +2.09%  '
$ws.Range('D19').Value = '''1.033'
$ws.Range('E19').Value = '  +2.60%  '
$ws.Range('D20').Value = '''15.50'
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('D21').Value = '27.749.13'
$ws.Range('E21').Value = '  +2.97%  '
$ws.Range('D22').Value = '''5.322'
$ws.Range('E22').Value = '  +2.42%  '
$ws.Range('D23').Value = '''11.46'
$ws.Range('E23').Value = '  +4.48%  '
$ws.Range('D24').Value = '''157.95'
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('D25').Value = '''1.943'
$ws.Range('E25').Value = '  +3.05%  '
$ws.Range('D26').Value = '''18.86'
$ws.Range('E26').Value = '  +2.92%  '
$ws.Range('D27').Value = '''1.995'
$ws.Range('E27').Value = '  +3.46%  '
$ws.Range('D28').Value = '''5.323'
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('D29').Value = '''117.55'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('D30').Value = '''0.09097'
$ws.Range('E30').Value = '  +1.60%  '
$ws.Range('D31').Value = '''1.217'
$ws.Range('E31').Value = '  +5.11%  '
$ws.Range('D32').Value = '''0.7695'
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('D33').Value = '''3.012'
$ws.Range('E33').Value = '  +7.31%  '
$ws.Range('D34').Value = '''4.581'
$ws.Range('E34').Value = '  +3.33%  '
$ws.Range('D35').Value = '''1.034'
$ws.Range('E35').Value = '  +2.64%  '
$ws.Range('D36').Value = '''1.163'
$ws.Range('E36').Value = '  +3.52%  '
$ws.Range('D37').Value = '''0.01992'
$ws.Range('E37').Value = '  +3.57%  '
$ws.Range('D38').Value = '''0.05356'
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('D39').Value = '''0.5204'
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('D40').Value = '''2.830'
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('D41').Value = '''0.1694'
$ws.Range('E41').Value = '  +2.52%  '
$ws.Range('D42').Value = '''6.854'
$ws.Range('E42').Value = '  +5.49%  '
$ws.Range('D43').Value = '''8.716'
$ws.Range('E43').Value = '  +4.75%  '
$ws.Range('D44').Value = '''109.61'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').Value = '''10.61'
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('D46').Value = '''1.731'
$ws.Range('E46').Value = '  +4.94%  '
$ws.Range('D47').Value = '''0.4699'
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('D48').Value = '''0.06424'
$ws.Range('E48').Value = '  +2.42%  '
$ws.Range('D49').Value = '''1.879'
$ws.Range('E49').Value = '  +3.76%  '
$ws.Range('D50').Value = '''39.85'
$ws.Range('E50').Value = '  +4.56%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''64.50'
$ws.Range('E51').Value = '  +1.48%  '
